$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Add the new "incomplete" worksheet right after Sheet1.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$newSheet.Name = "incomplete"

# Populate it with the (incomplete) employee data rows.
$newSheet.Range("A1").Value = "firstName"
$newSheet.Range("B1").Value = "middleName"
$newSheet.Range("C1").Value = "lastName"

$newSheet.Range("B2").Value = "MS"
$newSheet.Range("C2").Value = "data"

$newSheet.Range("A3").Value = "dummy"
$newSheet.Range("B3").Value = "MS"

# Match the cell selections captured in the saved workbook.
$newSheet.Range("E9").Select() | Out-Null

$sheet1.Activate() | Out-Null
$sheet1.Range("H12").Select() | Out-Null
